$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the full contents of row 2 and row 3 (columns A through AY).
$firstCol = 1    # A
$lastCol  = 51   # AY

# Columns whose populated values are text (not plain numbers/booleans).
# These need to be forced to "Text" number format before assignment so that
# numeric-looking strings (e.g. "1") and date-looking strings
# (e.g. "1988-01-01") are not reinterpreted by Excel as numbers/dates.
$textCols = @(3,4,6,7,8,9,12,16,20,21,22,23,25,26,27,28,29,35,49,50,51)

$row2Values = @()
$row3Values = @()

for ($col = $firstCol; $col -le $lastCol; $col++) {
    $row2Values += ,$ws.Cells.Item(2, $col).Value2
    $row3Values += ,$ws.Cells.Item(3, $col).Value2
}

for ($i = 0; $i -lt $row2Values.Length; $i++) {
    $col = $firstCol + $i
    $oldRow2Val = $row2Values[$i]
    $oldRow3Val = $row3Values[$i]

    # Nothing to do if the values already match between the two rows.
    if ($oldRow2Val -eq $oldRow3Val) {
        continue
    }

    $isText = $textCols -contains $col

    $cell2 = $ws.Cells.Item(2, $col)
    $cell3 = $ws.Cells.Item(3, $col)

    if ($isText) {
        $cell2.NumberFormat = "@"
        $cell3.NumberFormat = "@"
    }

    if ($null -eq $oldRow3Val) {
        $cell2.Value2 = ""
    } else {
        $cell2.Value2 = $oldRow3Val
    }

    if ($null -eq $oldRow2Val) {
        $cell3.Value2 = ""
    } else {
        $cell3.Value2 = $oldRow2Val
    }
}
